# Fixed problem with small files: update execution_datetime values for
# rows 9-13 and append a new row (14) for a newly processed small file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update execution_datetime for existing rows 9-13
$ws.Range("B9").Value  = 44375.90949464121
$ws.Range("B10").Value = 44375.90957033912
$ws.Range("B11").Value = 44375.90967922685
$ws.Range("B12").Value = 44375.90972654746
$ws.Range("B13").Value = 44375.90976891898

# Append the new row 14 for the small file that previously failed
$ws.Range("A14").Value = "ytube-transcripts-text---rZkdPXP6H4.txt"
$ws.Range("B14").Value = 44375.90801329167
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 34
